$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) and Column E (Volume 1h) updates ---
$ws.Range("D2").Value = "63.804.55"
$ws.Range("E2").Value = "  +0.69%  "

$ws.Range("D3").Value = "3.320.39"
$ws.Range("E3").Value = "  +2.57%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.24"
$ws.Range("E5").Value = "  +1.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.72"
$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.319.93"
$ws.Range("E8").Value = "  +2.71%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.55"
$ws.Range("E11").Value = "  +4.05%  "

$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.06"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").Value = "3.867.42"
$ws.Range("E15").Value = "  +2.61%  "

$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").Value = "3.319.36"
$ws.Range("E17").Value = "  +2.59%  "

$ws.Range("D18").Value = "63.902.33"
$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.77"
$ws.Range("E20").Value = "  +1.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.10"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("E22").Value = "  +2.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.98"
$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.96"
$ws.Range("E24").Value = "  +6.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.00"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("E27").Value = "  +1.85%  "

$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.25"
$ws.Range("E29").Value = "  +2.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("E30").Value = "  -4.21%  "

$ws.Range("E31").Value = "  +2.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.89"
$ws.Range("E32").Value = "  +5.24%  "

$ws.Range("E34").Value = "  +0.26%  "

$ws.Range("E35").Value = "  +1.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0399"
$ws.Range("E39").Value = "  +1.75%  "

$ws.Range("D40").Value = "3.133.01"
$ws.Range("E40").Value = "  +5.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "433.20"
$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("E42").Value = "  +7.05%  "

$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.75"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("E46").Value = "  +4.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.65"
$ws.Range("E47").Value = "  +8.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.40"
$ws.Range("E48").Value = "  +2.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("E51").Value = "  -0.57%  "

# --- Row 37/38: PEPE and OKB swap places with updated data ---
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.43"
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0745"
$ws.Range("E38").Value = "  +4.65%  "
